$d = $word.ActiveDocument

# 1. Merge the "Global Information Tracker" run and the closing curly-quote
#    run into a single run by replacing the combined text with itself; Word
#    collapses the match into one run using the formatting of the first run.
$d.Content.Find.Execute("Global Information Tracker”", $false, $false, $false, $false, $false, $true, 1, $false, "Global Information Tracker”", 2)

# 2. Insert a new run "Distribute VCS:- " right before the _GoBack bookmark
#    in the otherwise empty list paragraph.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Range.InsertBefore("Distribute VCS:- ")
